# "assume wkt is wgs84"
# - Drop the EPSG column (column E) from the AreaSource sheet: the workbook
#   no longer stores an EPSG code per area source, because the WKT geometry
#   is now always assumed to be WGS84 (lon/lat) rather than a projected CRS
#   such as EPSG:3857.
# - Rewrite the four existing WKT POLYGON strings (previously in a metric,
#   EPSG:3857-like projection) as WGS84 longitude/latitude coordinates.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AreaSource")

# Re-express the sample geometries in WGS84 lon/lat degrees instead of the
# old projected (metre-based) coordinates.
$ws.Range("D2").Value2 = "POLYGON ((22.006 41.444,22.016 41.443,22.015 41.438,22.009 41.437,22.005 41.438,22.006 41.444))"
$ws.Range("D3").Value2 = "POLYGON ((22.003 41.431,22.007 41.434,22.013 41.435,22.016 41.432,22.014 41.428,22.005 41.425,22.003 41.431))"
$ws.Range("D4").Value2 = "POLYGON ((21.9347884480412 41.4453304271493,21.9394977397296 41.4465438625105,21.9435006376646 41.4448891723006,21.9447662598059 41.4450656746002,21.9543614416207 41.4342319564415,21.9439421337604 41.4333934262775,21.9347884480412 41.4453304271493))"
$ws.Range("D5").Value2 = "POLYGON ((21.9645373511873 41.455713016914,21.9888354026551 41.4619862230537,21.987324040602 41.4423804391921,21.9665137477182 41.4438619720198,21.9645373511873 41.455713016914))"

# Drop the now-unnecessary EPSG column entirely (column E), shifting every
# later column one to the left.
$ws.Columns.Item(5).Delete()

# Move the active selection to where the deleted column used to be.
$ws.Range("E1").Select()
